# Weekly fruit/vegetable market data update.
# A new weekly record (row 32) is inserted above the existing "Espárragos"
# history, pushing the previous rows 32-43 down to 33-44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 32 (shifts old rows 32-43 -> 33-44,
# carrying their values/formatting down with them).
$ws.Rows(32).Insert()

# Populate the newly inserted row 32 with the new weekly observation.
$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(32, 3).Value = 'Ñuble'
$ws.Cells.Item(32, 4).Value = 44900
$ws.Cells.Item(32, 5).Value = 16
$ws.Cells.Item(32, 6).Value = 300000000
$ws.Cells.Item(32, 7).Value = 'Espárragos'
$ws.Cells.Item(32, 8).Value = 'Sin especificar'
$ws.Cells.Item(32, 9).Value = 'Primera'
$ws.Cells.Item(32, 10).Value = 1200
$ws.Cells.Item(32, 11).Value = 900
$ws.Cells.Item(32, 12).Value = 1000
$ws.Cells.Item(32, 13).Value = 950
$ws.Cells.Item(32, 14).Value = '$/kilo'
$ws.Cells.Item(32, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(32, 16).Value = 950
$ws.Cells.Item(32, 17).Value = 1
$ws.Cells.Item(32, 18).Value = 'Hortaliza'
